$d = $word.ActiveDocument

# 1. Title (paragraph 2)
$d.Paragraphs(2).Range.Text = 'Title: SXSW 2025 live coverage: The potential of quantum computing, Ireland''s prime minister makes a splash, and a Metallica concert in Apple Vision Pro'

# 2. Author(s) (paragraph 3)
$d.Paragraphs(3).Range.Text = 'Author(s): Sarah Perez, Kyle Wiggers, Kirsten Korosec, Rebecca Bellan, Lauren Forristal'

# 3. Published (paragraph 4)
$d.Paragraphs(4).Range.Text = 'Published: March 12, 2025 at 10:29 PM (EET)'

# 4. Summary body (paragraph 6)
$d.Paragraphs(6).Range.Text = 'SXSW 2025 kicked off Friday in Austin, Texas. The tech portion of the annual event will run through March 13. There are three $350,000 grants on offer and another for $125,000 to AI companies at any stage that are looking to create a presence in the city. Fogerty spoke at the SXSW conference in Austin, Texas. He talked about how the medium we use to consume music defines how it’s presented. He also touched on CCR history, music publishing rights, and the stories behind some of the band''s songs. Charina Chou, COO of Google Quantum AI, argued for government investment in the quantum computing market. Chou pointed out that China is currently outspending the U.S. two to one on quantum research. She also talked about the broader potential for quantum computing, the types of problems it could solve. Ireland’s Taoiseach (prime minister) Micheál Martin made a splash at SXSW with an onstage interview. His SXSW stop is part of a whirlwind tour that includes a visit to the White House. The tour comes at a critical time for Ireland, a small country of about 5 million people. Metallica frontman Lars Ulrich spoke at SXSW. He said there were 14 or 15 specialized cameras at the concert to capture the immersive experience. The songs from the concert will also be released in Spatial Audio on Apple Music. The city is launching a grant program to incentivize AI startups to set up shop in San Jose. San Jose is launching a grant program to incentivize AI startups to set up shop in the city. Conan O’Brien and Johanna Faries nerded out over the world of gaming at SXSW. IBM CEO Arvind Krishna said he thinks AI will ultimately make programmers more productive. IBM CEO Arvind Krishna says he thinks AI will ultimately make programmers more productive. Paramount CTO Phil Wiser remains optimistic about the future of AI. Wiser says he doesn’t think it can replace the role of a director or create a full-length movie. Filmmaking includes voice tools, dubbing, color correction, and de-aging. various other aspects of filmmaking, such as voice tools,. dubbing,. color correction,. de- aging. and other aspects. of filmmaking. For more information, visit www.filmmaking.com. for more information.'

Write-Output "replacements done"
